# Fruta / hortaliza, semanal
# Adds this week's record at the top (row 2), pushing the rest of the
# historical rows down by one, with the oldest row now appended at the
# bottom (row 30) with its original data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing data rows (2..29) down by one row, so that what used
# to be row 2 becomes row 3, ..., what used to be row 29 becomes row 30.
$ws.Rows("2:2").Insert()

# The inserted row inherits the header row's (bold/centered) formatting;
# reset it back to the plain style used by the rest of the data rows.
$ws.Range("A2:T2").Style = "Normal"

# Fill in the brand-new week's record in row 2.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 45083
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104003
$ws.Range("J2").Value = "Membrillo"
$ws.Range("K2").Value = "Champion"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 9000
$ws.Range("O2").Value = 10000
$ws.Range("P2").Value = 9500
$ws.Range("Q2").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R2").Value = "Región del Maule"
$ws.Range("S2").Value = 528
$ws.Range("T2").Value = 18
